# Update column C ("Förändrad") from 2023-11-13 (45243) to 2023-11-14 (45244)
# for all data rows (rows 2 through 21) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
